$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert the A1:K46 range into a native Excel Table (ListObject) ---
# This mirrors "Format as Table" in the Excel UI: the flat range becomes
# Table1 with headers taken from row 1, and the medium-blue banded style
# used in the target workbook.
$range = $ws.Range("A1:K46")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium9"

# --- Apply the AutoFilter so only entrp_ptnt_id = 1004 is shown ---
# Using the "discrete values" form (array + xlFilterValues) reproduces the
# <filters><filter val="1004"/></filters> shape (rather than a customFilter),
# and hides every row whose column-1 value isn't 1004 — matching the
# hidden="1" rows seen on everything except rows 22-24 (id 1004).
[void]$tbl.Range.AutoFilter(1, @("1004"), 7)

# --- Column widths (best effort; Excel quantizes these to its internal grid) ---
$ws.Columns.Item(1).ColumnWidth = 13.833333333333332
$ws.Columns.Item(2).ColumnWidth = 10.333333333333332
$ws.Columns.Item(3).ColumnWidth = 18.5
$ws.Columns.Item(5).ColumnWidth = 13.499999999999998
$ws.Columns.Item(6).ColumnWidth = 19.166666666666668
$ws.Columns.Item(7).ColumnWidth = 17.5
$ws.Columns.Item(8).ColumnWidth = 12.999999999999998
$ws.Columns.Item(9).ColumnWidth = 15.999999999999998
$ws.Columns.Item(10).ColumnWidth = 14.999999999999998
$ws.Columns.Item(11).ColumnWidth = 9.999999999999998

# --- Match the saved selection/active cell in the refreshed workbook ---
[void]$ws.Activate()
[void]$ws.Range("F8").Select()
